# Auto-generated Excel COM-interop edit script
# Updates cached numeric values (currentAveragePrice / LevePrice / LeveProfit columns)
# across the per-Job Leve tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to refreshed
# market-board figures, per the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 145.9
$ws.Range("I8").Value = 51
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 153
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = -14
$ws.Range("N8").Value = -3278
$ws.Range("H113").Value = 47623800
$ws.Range("J113").Value = 5812.5
$ws.Range("L113").Value = 5812.5
$ws.Range("N113").Value = -12320.5
$ws.Range("H129").Value = 233702.11
$ws.Range("J129").Value = 251208.53
$ws.Range("L129").Value = 753625.59
$ws.Range("N129").Value = -763625.59
$ws.Range("H132").Value = 3044.0908
$ws.Range("I132").Value = 3322.0386
$ws.Range("J132").Value = 2011.7142
$ws.Range("K132").Value = 9966.1158
$ws.Range("L132").Value = 6035.142599999999
$ws.Range("M132").Value = -7436.1158
$ws.Range("N132").Value = -11095.1426
$ws.Range("H135").Value = 11631245
$ws.Range("I135").Value = 689.0968
$ws.Range("K135").Value = 6201.8712
$ws.Range("M135").Value = -3666.8712
$ws.Range("H137").Value = 1808.5
$ws.Range("I137").Value = 1653.1786
$ws.Range("J137").Value = 2533.3333
$ws.Range("K137").Value = 4959.5358
$ws.Range("L137").Value = 7599.999899999999
$ws.Range("M137").Value = -2409.5358
$ws.Range("N137").Value = -12699.9999
$ws.Range("H138").Value = 10991226
$ws.Range("I138").Value = 22223078
$ws.Range("J138").Value = 3544.5652
$ws.Range("K138").Value = 66669234
$ws.Range("L138").Value = 10633.6956
$ws.Range("M138").Value = -66664094
$ws.Range("N138").Value = -20913.6956
$ws.Range("H141").Value = 1247.8
$ws.Range("I141").Value = 834.55817
$ws.Range("K141").Value = 2503.67451
$ws.Range("M141").Value = 2676.32549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 487722.56
$ws.Range("I61").Value = 621990.4399999999
$ws.Range("J61").Value = 1001.5
$ws.Range("K61").Value = 621990.4399999999
$ws.Range("L61").Value = 1001.5
$ws.Range("M61").Value = -621778.4399999999
$ws.Range("N61").Value = -1425.5
$ws.Range("H74").Value = 29413872
$ws.Range("I74").Value = 32260090
$ws.Range("K74").Value = 32260090
$ws.Range("M74").Value = -32259216
$ws.Range("H77").Value = 29413872
$ws.Range("I77").Value = 32260090
$ws.Range("K77").Value = 161300450
$ws.Range("M77").Value = -161296082
$ws.Range("H136").Value = 487722.56
$ws.Range("I136").Value = 621990.4399999999
$ws.Range("J136").Value = 1001.5
$ws.Range("K136").Value = 1865971.32
$ws.Range("L136").Value = 3004.5
$ws.Range("M136").Value = -1863421.32
$ws.Range("N136").Value = -8104.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 5500
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 10000
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = -860
$ws.Range("N8").Value = -10280
$ws.Range("H11").Value = 561.25
$ws.Range("I11").Value = 561.25
$ws.Range("K11").Value = 561.25
$ws.Range("M11").Value = -421.25
$ws.Range("H20").Value = 1701.8462
$ws.Range("I20").Value = 1891.5555
$ws.Range("J20").Value = 1275
$ws.Range("K20").Value = 1891.5555
$ws.Range("L20").Value = 1275
$ws.Range("M20").Value = -1644.5555
$ws.Range("N20").Value = -1769
$ws.Range("H22").Value = 289.65216
$ws.Range("I22").Value = 247.77777
$ws.Range("J22").Value = 440.4
$ws.Range("K22").Value = 247.77777
$ws.Range("L22").Value = 440.4
$ws.Range("M22").Value = -74.77777
$ws.Range("N22").Value = -786.4
$ws.Range("H86").Value = 1577.2307
$ws.Range("I86").Value = 1405
$ws.Range("J86").Value = 1778.1666
$ws.Range("K86").Value = 1405
$ws.Range("L86").Value = 1778.1666
$ws.Range("M86").Value = -282
$ws.Range("N86").Value = -4024.1666
$ws.Range("H89").Value = 1577.2307
$ws.Range("I89").Value = 1405
$ws.Range("J89").Value = 1778.1666
$ws.Range("K89").Value = 7025
$ws.Range("L89").Value = 8890.833000000001
$ws.Range("M89").Value = -1409
$ws.Range("N89").Value = -20122.833
$ws.Range("H134").Value = 3155.0889
$ws.Range("I134").Value = 3473.0286
$ws.Range("J134").Value = 2042.3
$ws.Range("K134").Value = 10419.0858
$ws.Range("L134").Value = 6126.9
$ws.Range("M134").Value = -7884.085800000001
$ws.Range("N134").Value = -11196.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 330.625
$ws.Range("I5").Value = 147.44444
$ws.Range("J5").Value = 440.53333
$ws.Range("K5").Value = 147.44444
$ws.Range("L5").Value = 440.53333
$ws.Range("M5").Value = -35.44443999999999
$ws.Range("N5").Value = -664.53333
$ws.Range("H22").Value = 333.8
$ws.Range("I22").Value = 292.25
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 292.25
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 57.75
$ws.Range("N22").Value = -1200
$ws.Range("H31").Value = 2904.1738
$ws.Range("I31").Value = 1669.0588
$ws.Range("K31").Value = 1669.0588
$ws.Range("M31").Value = -1374.0588
$ws.Range("H34").Value = 2904.1738
$ws.Range("I34").Value = 1669.0588
$ws.Range("K34").Value = 1669.0588
$ws.Range("M34").Value = -1467.0588
$ws.Range("H58").Value = 19788.666
$ws.Range("I58").Value = 1220
$ws.Range("J58").Value = 168338
$ws.Range("K58").Value = 1220
$ws.Range("L58").Value = 168338
$ws.Range("M58").Value = -1017
$ws.Range("N58").Value = -168744
$ws.Range("H132").Value = 1749.4897
$ws.Range("I132").Value = 1359.8085
$ws.Range("K132").Value = 4079.4255
$ws.Range("M132").Value = -1549.4255
$ws.Range("H136").Value = 19788.666
$ws.Range("I136").Value = 1220
$ws.Range("J136").Value = 168338
$ws.Range("K136").Value = 3660
$ws.Range("L136").Value = 505014
$ws.Range("M136").Value = -1110
$ws.Range("N136").Value = -510114

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 10000
$ws.Range("J74").Value = 10000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -32122
$ws.Range("H77").Value = 10000
$ws.Range("J77").Value = 10000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -100608
$ws.Range("H131").Value = 728.86
$ws.Range("J131").Value = 738.9053
$ws.Range("L131").Value = 2216.7159
$ws.Range("N131").Value = -12296.7159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 68.1875
$ws.Range("I2").Value = 63.1
$ws.Range("J2").Value = 76.666664
$ws.Range("K2").Value = 63.1
$ws.Range("L2").Value = 76.666664
$ws.Range("M2").Value = 49.9
$ws.Range("N2").Value = -302.666664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 252.8
$ws.Range("I16").Value = 274.76923
$ws.Range("J16").Value = 212
$ws.Range("K16").Value = 274.76923
$ws.Range("L16").Value = 212
$ws.Range("M16").Value = -104.76923
$ws.Range("N16").Value = -552
$ws.Range("H22").Value = 6750.3335
$ws.Range("I22").Value = 5125.5
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 5125.5
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -4830.5
$ws.Range("N22").Value = -10590
$ws.Range("H27").Value = 6750.3335
$ws.Range("I27").Value = 5125.5
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 5125.5
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = -5018.5
$ws.Range("N27").Value = -10214
$ws.Range("H55").Value = 162.125
$ws.Range("I55").Value = 162.70589
$ws.Range("J55").Value = 160.71428
$ws.Range("K55").Value = 162.70589
$ws.Range("L55").Value = 160.71428
$ws.Range("M55").Value = 10.29410999999999
$ws.Range("N55").Value = -506.71428
$ws.Range("H132").Value = 1038.7115
$ws.Range("I132").Value = 1049.2745
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 3147.8235
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -617.8235
$ws.Range("N132").Value = -6560
$ws.Range("H136").Value = 1506.3334
$ws.Range("I136").Value = 1433.6
$ws.Range("J136").Value = 1870
$ws.Range("K136").Value = 4300.799999999999
$ws.Range("L136").Value = 5610
$ws.Range("M136").Value = -1750.799999999999
$ws.Range("N136").Value = -10710

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1375.4445
$ws.Range("I96").Value = 1250
$ws.Range("J96").Value = 1532.25
$ws.Range("K96").Value = 1250
$ws.Range("L96").Value = 1532.25
$ws.Range("M96").Value = 123
$ws.Range("N96").Value = -4278.25
$ws.Range("H132").Value = 534.8837
$ws.Range("I132").Value = 546.1905
$ws.Range("K132").Value = 1638.5715
$ws.Range("M132").Value = 891.4285
$ws.Range("H136").Value = 17243474
$ws.Range("I136").Value = 24391222
$ws.Range("K136").Value = 73173666
$ws.Range("M136").Value = -73171116
